$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of fruit/vegetable price data: reshuffle Fecha, Calidad, Volumen,
# Precio minimo/maximo/promedio, Origen and Precio $/Kg across the Achicoria rows.

# Row 2
$ws.Range("D2").Value = 44389
$ws.Range("J2").Value = 55

# Row 3
$ws.Range("D3").Value = 44369
$ws.Range("J3").Value = 60

# Row 4
$ws.Range("D4").Value = 44420
$ws.Range("J4").Value = 45

# Row 5
$ws.Range("D5").Value = 44355
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 8000
$ws.Range("M5").Value = 8000
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 500

# Row 6
$ws.Range("D6").Value = 44467
$ws.Range("J6").Value = 40
$ws.Range("O6").Value = "Región del Maule"

# Row 7
$ws.Range("D7").Value = 44348
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 35
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 7000
$ws.Range("P7").Value = 438

# Row 8
$ws.Range("D8").Value = 44305
$ws.Range("J8").Value = 35

# Row 9
$ws.Range("D9").Value = 44397
$ws.Range("J9").Value = 40

# Row 10
$ws.Range("D10").Value = 44308
$ws.Range("J10").Value = 75
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 5000
$ws.Range("M10").Value = 5000
$ws.Range("P10").Value = 312

# Row 11
$ws.Range("D11").Value = 44398
$ws.Range("J11").Value = 80
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("P11").Value = 438

# Row 12
$ws.Range("D12").Value = 44362
$ws.Range("J12").Value = 25
$ws.Range("K12").Value = 8000
$ws.Range("L12").Value = 8000
$ws.Range("M12").Value = 8000
$ws.Range("P12").Value = 500

# Row 13
$ws.Range("D13").Value = 44313
$ws.Range("J13").Value = 20
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = 7000
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 438

# Row 14
$ws.Range("D14").Value = 44403
$ws.Range("J14").Value = 35
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = 5000
$ws.Range("P14").Value = 312

# Row 15
$ws.Range("D15").Value = 44386
$ws.Range("J15").Value = 40
$ws.Range("O15").Value = "Región del Maule"

# Row 16
$ws.Range("D16").Value = 44371
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 7000
$ws.Range("P16").Value = 438

# Row 17
$ws.Range("D17").Value = 44396
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 7000
$ws.Range("M17").Value = 7000
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 438

# Row 18
$ws.Range("D18").Value = 44354
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 8000
$ws.Range("L18").Value = 9000
$ws.Range("M18").Value = 8500
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 531

# Row 19
$ws.Range("D19").Value = 44354
$ws.Range("J19").Value = 80
$ws.Range("K19").Value = 9000
$ws.Range("L19").Value = 9000
$ws.Range("M19").Value = 9000
$ws.Range("P19").Value = 562

# Row 20
$ws.Range("D20").Value = 44399
$ws.Range("J20").Value = 80
$ws.Range("O20").Value = "Región Metropolitana"

# Row 21
$ws.Range("D21").Value = 44372
$ws.Range("J21").Value = 50
$ws.Range("K21").Value = 6000
$ws.Range("M21").Value = 6400
$ws.Range("P21").Value = 400

# Row 22
$ws.Range("D22").Value = 44315
$ws.Range("J22").Value = 40

# Row 23
$ws.Range("D23").Value = 44312
$ws.Range("J23").Value = 40

# Row 24
$ws.Range("D24").Value = 44314
$ws.Range("I24").Value = "Segunda"
$ws.Range("J24").Value = 20
$ws.Range("K24").Value = 5000
$ws.Range("L24").Value = 5000
$ws.Range("M24").Value = 5000
$ws.Range("O24").Value = "Región del Maule"
$ws.Range("P24").Value = 312

# Row 25
$ws.Range("D25").Value = 44385
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 7000
$ws.Range("L25").Value = 7000
$ws.Range("M25").Value = 7000
$ws.Range("P25").Value = 438

# Row 26
$ws.Range("D26").Value = 44392
$ws.Range("J26").Value = 95
